$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rsquo = [char]0x2019
$mdash = [char]0x2014

# Delete rows 8-35 (the old per-field rows no longer needed)
$ws.Range("A8:A35").EntireRow.Delete() | Out-Null

# Set consolidated tuple-like strings for rows 2-7
$ws.Range("A2").Value = "('Breaking', ['{U}{B}', 'Sorcery', 'Target player mills eight cards.', 'Fuse (You may cast one or both halves of this card from your hand.)', 'Entering', '{4}{B}{R}', 'Sorcery', 'Put a creature card from a graveyard onto the battlefield under your control. It gains haste until end of turn.', 'Fuse (You may cast one or both halves of this card from your hand.)'])"

$ws.Range("A3").Value = "(`"Maze's End`", ['Land', 'Maze${rsquo}s End enters the battlefield tapped.', '{T}: Add {C}.', '{3}, {T}, Return Maze${rsquo}s End to its owner${rsquo}s hand: Search your library for a Gate card, put it onto the battlefield, then shuffle your library. If you control ten or more Gates with different names, you win the game.'])"

$ws.Range("A4").Value = "('Melek, Izzet Paragon', ['{4}{U}{R}', 'Legendary Creature ${mdash} Weird Wizard', 'Play with the top card of your library revealed.', 'You may cast instant and sorcery spells from the top of your library.', 'Whenever you cast an instant or sorcery spell from your library, copy it. You may choose new targets for the copy.', '2/4'])"

$ws.Range("A5").Value = "('Plains', ['Basic Land ${mdash} Plains', '({T}: Add {W}.)'])"

$ws.Range("A6").Value = "('Render Silent', ['{W}{U}{U}', 'Instant', 'Counter target spell. Its controller can${rsquo}t cast spells this turn.'])"

$ws.Range("A7").Value = "(`"Trostani's Summoner`", ['{5}{G}{W}', 'Creature ${mdash} Elf Shaman', 'When Trostani${rsquo}s Summoner enters the battlefield, create a 2/2 white Knight creature token with vigilance, a 3/3 green Centaur creature token, and a 4/4 green Rhino creature token with trample.', '1/1'])"
